# The presentation ships with two DrawingML theme parts:
#   ppt/theme/theme1.xml  -> originally "Office Theme" (used by the Notes Master)
#   ppt/theme/theme2.xml  -> originally "Integral"      (used by the Slide Master /
#                                                         the presentation's live design)
#
# The authored change swaps the two themes' contents wholesale: theme1.xml becomes
# the "Integral" theme and theme2.xml becomes the "Office Theme". In other words,
# the presentation's active/applied design theme (the one driving every slide via
# the Slide Master) switches from "Integral" colours to the stock "Office Theme"
# colours.
#
# The PowerPoint object model's ThemeColorScheme maps 1:1 onto the <a:clrScheme>
# children (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) of the theme that is
# wired to the Slide Master (theme2.xml), so we drive the swap through it.

function Convert-HexToComColor($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

# Target values = the stock "Office Theme" colour scheme (what used to live in
# theme1.xml before the swap).
$officeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $officeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = Convert-HexToComColor $officeColors[$i - 1]
}
